# Update "想去人数" (F column, interest counts) and a couple of "最低票价" (G column)
# values across sheets 展览, 演出, and 全部类型, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 526
$ws.Range("F9").Value = 1005
$ws.Range("F11").Value = 230
$ws.Range("F14").Value = 807
$ws.Range("F18").Value = 1322
$ws.Range("F19").Value = 118
$ws.Range("F20").Value = 844
$ws.Range("F21").Value = 1159
$ws.Range("F23").Value = 1378
$ws.Range("F24").Value = 685
$ws.Range("F25").Value = 186
$ws.Range("F26").Value = 1264
$ws.Range("F29").Value = 345
$ws.Range("F30").Value = 2901
$ws.Range("F31").Value = 566
$ws.Range("F33").Value = 1380

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G5").Value = 80
$ws.Range("F7").Value = 10
$ws.Range("F10").Value = 152

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 526
$ws.Range("G10").Value = 80
$ws.Range("F14").Value = 10
$ws.Range("F16").Value = 1005
$ws.Range("F18").Value = 230
$ws.Range("F22").Value = 152
$ws.Range("F26").Value = 807
$ws.Range("F30").Value = 1322
$ws.Range("F31").Value = 118
$ws.Range("F32").Value = 844
$ws.Range("F33").Value = 1159
$ws.Range("F35").Value = 1378
$ws.Range("F36").Value = 685
$ws.Range("F37").Value = 186
$ws.Range("F38").Value = 1264
$ws.Range("F43").Value = 345
$ws.Range("F44").Value = 2901
$ws.Range("F45").Value = 566
$ws.Range("F47").Value = 1380
